# Updated cryptos list refresh (Coin / Link / Price / Volume(1h) columns on
# Sheet1) to match the latest scrape.
#
# The sheet stores every cell as text (prices use a European "."-grouped
# display, e.g. "71.872.10", and volumes are padded percentage strings like
# "  -2.05%  "). Plain decimal-looking price updates (e.g. "598.58") would
# normally be auto-coerced to a number by Excel on assignment, so those are
# written with a leading apostrophe (forcing text entry, same as typing
# '598.58 into a cell) and the style is reset back to "Normal" afterwards so
# no stray quote-prefix/number-format styling is left on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "71.872.10" },
    @{ Cell = "E2"; Value = "  -2.05%  " },
    @{ Cell = "D3"; Value = "2.662.70" },
    @{ Cell = "E3"; Value = "  +0.00%  " },
    @{ Cell = "E4"; Value = "  -0.04%  " },
    @{ Cell = "D5"; Value = "598.58" },
    @{ Cell = "E5"; Value = "  -1.56%  " },
    @{ Cell = "D6"; Value = "173.93" },
    @{ Cell = "E6"; Value = "  -3.55%  " },
    @{ Cell = "E7"; Value = "  +0.01%  " },
    @{ Cell = "D8"; Value = "0.524" },
    @{ Cell = "E8"; Value = "  -1.20%  " },
    @{ Cell = "D9"; Value = "2.660.49" },
    @{ Cell = "E9"; Value = "  -0.07%  " },
    @{ Cell = "E10"; Value = "  -4.07%  " },
    @{ Cell = "E11"; Value = "  +2.20%  " },
    @{ Cell = "D12"; Value = "0.355" },
    @{ Cell = "E12"; Value = "  +0.11%  " },
    @{ Cell = "E13"; Value = "  -2.16%  " },
    @{ Cell = "D14"; Value = "3.154.77" },
    @{ Cell = "E14"; Value = "  +0.54%  " },
    @{ Cell = "D15"; Value = "0.0000184" },
    @{ Cell = "E15"; Value = "  -3.97%  " },
    @{ Cell = "D16"; Value = "71.823.92" },
    @{ Cell = "E16"; Value = "  -1.95%  " },
    @{ Cell = "D17"; Value = "26.21" },
    @{ Cell = "E17"; Value = "  -2.93%  " },
    @{ Cell = "D18"; Value = "2.668.90" },
    @{ Cell = "E18"; Value = "  +1.04%  " },
    @{ Cell = "D19"; Value = "12.20" },
    @{ Cell = "E19"; Value = "  +4.82%  " },
    @{ Cell = "E20"; Value = "  +1.74%  " },
    @{ Cell = "D21"; Value = "371.91" },
    @{ Cell = "E21"; Value = "  -4.07%  " },
    @{ Cell = "E22"; Value = "  -1.49%  " },
    @{ Cell = "D23"; Value = "2.03" },
    @{ Cell = "E23"; Value = "  -0.27%  " },
    @{ Cell = "D24"; Value = "72.01" },
    @{ Cell = "E24"; Value = "  -2.17%  " },
    @{ Cell = "E25"; Value = "  -0.06%  " },
    @{ Cell = "D26"; Value = "4.33" },
    @{ Cell = "E26"; Value = "  -2.54%  " },
    @{ Cell = "D27"; Value = "9.75" },
    @{ Cell = "E27"; Value = "  -2.04%  " },
    @{ Cell = "D28"; Value = "2.803.11" },
    @{ Cell = "E28"; Value = "  +0.00%  " },
    @{ Cell = "D29"; Value = "0.998" },
    @{ Cell = "D30"; Value = "0.0₃0967" },
    @{ Cell = "E30"; Value = "  -0.90%  " },
    @{ Cell = "E31"; Value = "  -0.68%  " },
    @{ Cell = "D32"; Value = "501.05" },
    @{ Cell = "E32"; Value = "  -6.83%  " },
    @{ Cell = "E33"; Value = "  -3.74%  " },
    @{ Cell = "E34"; Value = "  -1.31%  " },
    @{ Cell = "E35"; Value = "  -0.09%  " },
    @{ Cell = "D36"; Value = "163.04" },
    @{ Cell = "E36"; Value = "  -0.68%  " },
    @{ Cell = "D37"; Value = "19.51" },
    @{ Cell = "E37"; Value = "  +0.44%  " },
    @{ Cell = "E38"; Value = "  -0.38%  " },
    @{ Cell = "E39"; Value = "  -3.37%  " },
    @{ Cell = "E40"; Value = "  -3.31%  " },
    @{ Cell = "D41"; Value = "1.76" },
    @{ Cell = "E41"; Value = "  -4.96%  " },
    @{ Cell = "E42"; Value = "  -0.10%  " },
    @{ Cell = "E43"; Value = "  -3.45%  " },
    @{ Cell = "E44"; Value = "  -3.74%  " },
    @{ Cell = "E45"; Value = "  -1.13%  " },
    @{ Cell = "B46"; Value = "Aave" },
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" },
    @{ Cell = "D46"; Value = "156.28" },
    @{ Cell = "E46"; Value = "  +2.34%  " },
    @{ Cell = "B47"; Value = "OKB" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb" },
    @{ Cell = "D47"; Value = "39.45" },
    @{ Cell = "E47"; Value = "  -0.97%  " },
    @{ Cell = "D48"; Value = "0.559" },
    @{ Cell = "E48"; Value = "  +2.21%  " },
    @{ Cell = "E49"; Value = "  +0.22%  " },
    @{ Cell = "E50"; Value = "  +0.77%  " },
    @{ Cell = "D51"; Value = "0.0755" },
    @{ Cell = "E51"; Value = "  -1.79%  " }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $value = $update.Value

    if ($value -match '^[+-]?\d+(\.\d+)?$') {
        # Looks like a plain number - force text so Excel keeps the original
        # textual representation instead of parsing it into a numeric value.
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
